$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - mirror style of existing header cells (H1 etc.)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows for column I (I0) and J (IF)
$data = @(
    @(1, 3),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(6, 9),
    @(2, 3),
    @(6, 6),
    @(8, 8),
    @(4, 5),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$wb.Save()
